$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update header summary values ---
# "VALOR MORA" total
$ws.Range("E11").Value = 40674
# "Cant. Trabajadores" / "Cant. Periodos"
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

# --- Rebuild the detail table (rows 16-20) down to just 2 workers ---
# Row 17 (RICARDO's current row) needs to inherit the "last row" bottom-border
# formatting that currently belongs to row 20 (BRAHIAN), since after the
# update the table only has 2 data rows and row 17 becomes the last one.
$ws.Range("B20:J20").Copy() | Out-Null
$ws.Range("B17:J17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Clear out the old worker rows' content first so their text values become
# unused and the database can be repopulated fresh.
$ws.Range("B16:J20").ClearContents() | Out-Null

# New row 16: RAYMUNDO JOSE FLOREZ SANTOYA
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047492513"
$ws.Range("D16").Value = "RAYMUNDO JOSE FLOREZ SANTOYA"
$ws.Range("E16").Value = "2505"
$ws.Range("F16").Value = 16000
$ws.Range("G16").Value = 4000000

# New row 17: RICARDO ANTONIO MACHACON FAJARDO
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1143324982"
$ws.Range("D17").Value = "RICARDO ANTONIO MACHACON FAJARDO"
$ws.Range("E17").Value = "2506"
$ws.Range("F17").Value = 24674
$ws.Range("G17").Value = 1423500

# Remove the now-obsolete rows (old rows 18, 19, 20); this also shifts the
# signature block (previously rows 25-26) up to rows 22-23.
$ws.Rows("18:20").Delete() | Out-Null
